$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 2, pushing the existing
# data (old rows 2-9) down to rows 4-11.
$insertRange = $ws.Range("A2:F3")
$insertRange.EntireRow.Insert()

# Fill in the two newly inserted rows.
$ws.Range("A2").Value = "Norfolk"
$ws.Range("B2").Value = "NE"
$ws.Range("C2").Value = "Oglesby"
$ws.Range("C2").ClearFormats()
$ws.Range("C2").Font.Size = 11
$ws.Range("D2").Value = "IL"
$ws.Range("E2").Value = "<50000.0"
$ws.Range("F2").Value = "Skip"

$ws.Range("A3").Value = "Norfolk"
$ws.Range("B3").Value = "NE"
$ws.Range("C3").Value = "Any"
$ws.Range("D3").Value = "KS"
$ws.Range("E3").Value = "<50000.0"
$ws.Range("F3").Value = "Skip"

# Update the selection to match the edited workbook.
$ws.Range("G2").Select()
